$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J53").Value = 1486.125
$ws.Range("L53").Value = 1486.125
$ws.Range("N53").Value = -2760.125

$ws.Range("H55").Value = 937.2
$ws.Range("I55").Value = 240.3
$ws.Range("J55").Value = 2331
$ws.Range("K55").Value = 240.3
$ws.Range("L55").Value = 2331
$ws.Range("M55").Value = -26.30000000000001
$ws.Range("N55").Value = -2759

$ws.Range("H88").Value = 1970.25
$ws.Range("J88").Value = 2344.5
$ws.Range("L88").Value = 2344.5
$ws.Range("N88").Value = -3156.5

$ws.Range("H91").Value = 1970.25
$ws.Range("J91").Value = 2344.5
$ws.Range("L91").Value = 2344.5
$ws.Range("N91").Value = -5152.5

$ws.Range("H118").Value = 677.8
$ws.Range("I118").Value = 677.8
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2033.4
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -376.3999999999999
$ws.Range("N118").ClearContents()

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H138").Value = 3675.9297
$ws.Range("I138").Value = 2286
$ws.Range("J138").Value = 4635.643
$ws.Range("K138").Value = 6858
$ws.Range("L138").Value = 13906.929
$ws.Range("M138").Value = -1718
$ws.Range("N138").Value = -24186.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7622.4795
$ws.Range("I32").Value = 5382.0312
$ws.Range("J32").Value = 23554.555
$ws.Range("K32").Value = 5382.0312
$ws.Range("L32").Value = 23554.555
$ws.Range("M32").Value = -5095.0312
$ws.Range("N32").Value = -24128.555

$ws.Range("H132").Value = 1767.75
$ws.Range("I132").Value = 1827.025
$ws.Range("J132").Value = 1175
$ws.Range("K132").Value = 5481.075000000001
$ws.Range("L132").Value = 3525
$ws.Range("M132").Value = -2951.075000000001
$ws.Range("N132").Value = -8585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 7533.3335
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20504

$ws.Range("H36").Value = 2333
$ws.Range("I36").Value = 2333
$ws.Range("K36").Value = 2333
$ws.Range("M36").Value = -1799

$ws.Range("H37").Value = 5668.25
$ws.Range("I37").Value = 4224.3335
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 4224.3335
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -4087.3335
$ws.Range("N37").Value = -10274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27031542
$ws.Range("I31").Value = 52634812
$ws.Range("K31").Value = 52634812
$ws.Range("M31").Value = -52634517

$ws.Range("H34").Value = 27031542
$ws.Range("I34").Value = 52634812
$ws.Range("K34").Value = 52634812
$ws.Range("M34").Value = -52634610

$ws.Range("H58").Value = 1808.6666
$ws.Range("I58").Value = 1528.3784
$ws.Range("J58").Value = 3105
$ws.Range("K58").Value = 1528.3784
$ws.Range("L58").Value = 3105
$ws.Range("M58").Value = -1325.3784
$ws.Range("N58").Value = -3511

$ws.Range("H132").Value = 1814.6177
$ws.Range("I132").Value = 1714.9
$ws.Range("J132").Value = 2562.5
$ws.Range("K132").Value = 5144.700000000001
$ws.Range("L132").Value = 7687.5
$ws.Range("M132").Value = -2614.700000000001
$ws.Range("N132").Value = -12747.5

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 1433.3513
$ws.Range("I134").Value = 1121.6552
$ws.Range("J134").Value = 2563.25
$ws.Range("K134").Value = 3364.9656
$ws.Range("L134").Value = 7689.75
$ws.Range("M134").Value = -829.9655999999995
$ws.Range("N134").Value = -12759.75

$ws.Range("H135").Value = 130000
$ws.Range("J135").Value = 130000
$ws.Range("L135").Value = 130000
$ws.Range("N135").Value = -140140

$ws.Range("H136").Value = 1808.6666
$ws.Range("I136").Value = 1528.3784
$ws.Range("J136").Value = 3105
$ws.Range("K136").Value = 4585.135200000001
$ws.Range("L136").Value = 9315
$ws.Range("M136").Value = -2035.135200000001
$ws.Range("N136").Value = -14415

$ws.Range("H137").Value = 99129.5
$ws.Range("J137").Value = 99129.5
$ws.Range("L137").Value = 99129.5
$ws.Range("N137").Value = -109329.5

$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws.Range("H141").Value = 488716.56
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 488716.56
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 488716.56
$ws.Range("N141").Value = -499076.56
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 104.44444
$ws.Range("I2").Value = 109.75
$ws.Range("J2").Value = 93.833336
$ws.Range("K2").Value = 658.5
$ws.Range("L2").Value = 563.000016
$ws.Range("M2").Value = -545.5
$ws.Range("N2").Value = -789.000016

$ws.Range("H5").Value = 746.2353
$ws.Range("I5").Value = 705.8
$ws.Range("J5").Value = 804
$ws.Range("K5").Value = 2117.4
$ws.Range("L5").Value = 2412
$ws.Range("M5").Value = -2005.4
$ws.Range("N5").Value = -2636

$ws.Range("H14").Value = 19141.818
$ws.Range("I14").Value = 19141.818
$ws.Range("K14").Value = 57425.454
$ws.Range("M14").Value = -57252.454

$ws.Range("H33").Value = 4620406
$ws.Range("I33").Value = 206.42857
$ws.Range("J33").Value = 8663081
$ws.Range("K33").Value = 1238.57142
$ws.Range("L33").Value = 51978486
$ws.Range("M33").Value = -955.57142
$ws.Range("N33").Value = -51979052

$ws.Range("H38").Value = 177.5
$ws.Range("I38").Value = 20
$ws.Range("J38").Value = 230
$ws.Range("K38").Value = 60
$ws.Range("L38").Value = 690
$ws.Range("M38").Value = 287
$ws.Range("N38").Value = -1384

$ws.Range("H56").Value = 14824.16
$ws.Range("I56").Value = 14824.16
$ws.Range("K56").Value = 14824.16
$ws.Range("M56").Value = -14294.16

$ws.Range("H123").Value = 12033.25
$ws.Range("J123").Value = 33333
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -104899

$ws.Range("H124").Value = 2300
$ws.Range("I124").Value = 2300
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 6900
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -1990
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 27719
$ws.Range("I125").Value = 19299
$ws.Range("J125").Value = 33332.332
$ws.Range("K125").Value = 57897
$ws.Range("L125").Value = 99996.99600000001
$ws.Range("M125").Value = -52977
$ws.Range("N125").Value = -109836.996

$ws.Range("H129").Value = 11910760
$ws.Range("I129").Value = 41672456
$ws.Range("J129").Value = 6083
$ws.Range("K129").Value = 125017368
$ws.Range("L129").Value = 18249
$ws.Range("M129").Value = -125012368
$ws.Range("N129").Value = -28249

$ws.Range("H130").Value = 14426.2
$ws.Range("I130").Value = 4899.5
$ws.Range("J130").Value = 20777.334
$ws.Range("K130").Value = 14698.5
$ws.Range("L130").Value = 62332.00199999999
$ws.Range("M130").Value = -9678.5
$ws.Range("N130").Value = -72372.002

$ws.Range("H135").Value = 746.2353
$ws.Range("I135").Value = 705.8
$ws.Range("J135").Value = 804
$ws.Range("K135").Value = 6352.2
$ws.Range("L135").Value = 7236
$ws.Range("M135").Value = -3817.2
$ws.Range("N135").Value = -12306

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3648.111
$ws.Range("I22").Value = 3639.1667
$ws.Range("J22").Value = 3666
$ws.Range("K22").Value = 3639.1667
$ws.Range("L22").Value = 3666
$ws.Range("M22").Value = -3344.1667
$ws.Range("N22").Value = -4256

$ws.Range("H27").Value = 3648.111
$ws.Range("I27").Value = 3639.1667
$ws.Range("J27").Value = 3666
$ws.Range("K27").Value = 3639.1667
$ws.Range("L27").Value = 3666
$ws.Range("M27").Value = -3532.1667
$ws.Range("N27").Value = -3880

$ws.Range("H132").Value = 2684.95
$ws.Range("I132").Value = 1726.8695
$ws.Range("K132").Value = 5180.6085
$ws.Range("M132").Value = -2650.6085
